$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.731.75"
$ws.Range("E2").Value = "  +2.65%  "

$ws.Range("D3").Value = "1.789.45"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'223.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "

$ws.Range("D6").Value = "'0.555"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.61%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'32.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.40%  "

$ws.Range("E9").Value = "  +0.56%  "

$ws.Range("D10").Value = "'0.0688"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.49%  "

$ws.Range("D11").Value = "'0.0936"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.50%  "

$ws.Range("D12").Value = "2.045.85"
$ws.Range("E12").Value = "  +0.58%  "

$ws.Range("D13").Value = "'11.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.75%  "

$ws.Range("D14").Value = "1.795.62"
$ws.Range("E14").Value = "  +1.03%  "

$ws.Range("D15").Value = "34.725.13"
$ws.Range("E15").Value = "  +2.73%  "

$ws.Range("E16").Value = "  +1.19%  "

$ws.Range("E17").Value = "  +3.35%  "

$ws.Range("D18").Value = "'68.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").Value = "'253.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.00%  "

$ws.Range("D20").Value = "0.0₃0784"
$ws.Range("E20").Value = "  +6.10%  "

$ws.Range("D21").Value = "'1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").Value = "'10.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.06%  "

$ws.Range("E23").Value = "  +0.69%  "

$ws.Range("E24").Value = "  -0.82%  "

$ws.Range("D25").Value = "'159.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.58%  "

$ws.Range("D26").Value = "'16.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.34%  "

$ws.Range("E27").Value = "  +1.47%  "

$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  -1.35%  "

$ws.Range("E31").Value = "  +0.09%  "

$ws.Range("E32").Value = "  -0.37%  "

$ws.Range("D33").Value = "'3.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("E34").Value = "  +0.80%  "

$ws.Range("D35").Value = "1.438.08"
$ws.Range("E35").Value = "  -2.67%  "

$ws.Range("E36").Value = "  -1.34%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.0189"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.46%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.631"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.28%  "

$ws.Range("D39").Value = "'82.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("D40").Value = "'2.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.15%  "

$ws.Range("D42").Value = "'0.905"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.17%  "

$ws.Range("E43").Value = "  -0.53%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'1.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").Value = "'0.0504"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.69%  "

$ws.Range("E46").Value = "  +4.33%  "

$ws.Range("D47").Value = "1.943.44"
$ws.Range("E47").Value = "  +0.63%  "

$ws.Range("D48").Value = "'105.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.84%  "

$ws.Range("D49").Value = "'11.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.15%  "

$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("D51").Value = "'49.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.37%  "
